$d = $word.ActiveDocument

# --- 1. Fix the heading wording -----------------------------------------
# The heading currently reads "Objetivos específicos son S.M.A.R.T" which
# has a spurious "son" in it; the correction drops that word so it reads
# "Objetivos específicos S.M.A.R.T".
$heading = $d.Paragraphs(1).Range

# Word always keeps a "_GoBack" bookmark marking the last edited spot, and
# moves it whenever a new edit is made. Drop the stale one (it currently
# sits in the trailing empty paragraph) before editing.
try {
    $d.Bookmarks("_GoBack").Delete()
} catch {
    # no pre-existing _GoBack bookmark - nothing to remove
}

$headingFind = $heading.Find
$headingFind.ClearFormatting()
$headingFind.Execute("específicos son S.M.A.R.T", $true, $false, $false, $false, $false, `
                      $true, 1, $false, "específicos S.M.A.R.T", 2)

# --- 2. Re-drop "_GoBack" at the point of the edit ----------------------
# Word stamps the bookmark right where the cursor ends up after the edit,
# i.e. right before "S.M.A.R.T".
$headingRange = $d.Paragraphs(1).Range
$headingText = $headingRange.Text
$splitOffset = $headingRange.Start + $headingText.IndexOf("S.M.A.R.T")
$editPoint = $d.Range($splitOffset, $splitOffset)
$d.Bookmarks.Add("_GoBack", $editPoint)
